# ---------------------------------------------------------------------------
# Feb 6.docx lab-write-up edit:
#   1. "CHECK(ACTNO>100)"  ->  "CHECK(ACTNO>=100)"   (insert "=" just after
#      the "&gt;" and drop the editing-cursor "_GoBack" bookmark there,
#      which is where Word last left it when the file was saved)
#   2. the old "_GoBack" bookmark that used to sit between
#      "SELECT * FROM DEPOSIT_21052264" and " ORDER BY AMOUNT DESC;" is gone
#      (cursor moved away) so those two runs collapse back into one run.
#   3. table style "Normal Table" (styleId 3) picks up <w:qFormat/>, and its
#      matching <w:lsdException .../> in the latent-style table gains
#      w:qFormat="1".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. ALTER TABLE ... CHECK(ACTNO>100) -> CHECK(ACTNO>=100) --------------

$chk = $d.Content.Duplicate
$chk.Find.Execute("CHECK(ACTNO>") | Out-Null
$chk.Collapse(0)
$chk.InsertAfter("=")

# Collapsed range right after the newly-typed "=" is where Word stamps the
# "last edit" _GoBack bookmark; Bookmarks.Add with the reserved name
# "_GoBack" both (re)creates it there and removes it from its previous spot.
$eq = $d.Content.Duplicate
$eq.Find.Execute("CHECK(ACTNO>=") | Out-Null
$eq.Collapse(0)
$d.Bookmarks.Add("_GoBack", $eq)

# --- 2. SELECT * FROM DEPOSIT_21052264 [_GoBack] ORDER BY AMOUNT DESC; -----
# The bookmark move above already deleted the old _GoBack that used to split
# this sentence into two runs; re-assert the sentence as a single piece of
# text so the two runs coalesce into one.

$d.Content.Find.Execute("SELECT * FROM DEPOSIT_21052264 ORDER BY AMOUNT DESC;", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "SELECT * FROM DEPOSIT_21052264 ORDER BY AMOUNT DESC;", 2) | Out-Null

# --- 3. styles.xml: "Normal Table" becomes qFormat -------------------------

foreach ($sty in $d.Styles) {
    if ($sty.NameLocal -eq "Normal Table") {
        $sty.QuickStyle = $true
    }
}
